$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 6 (G=4564)
$ws.Range("H6").Value = 290.7143
$ws.Range("I6").Value = 23.333334
$ws.Range("J6").Value = 772
$ws.Range("K6").Value = 70.00000199999999
$ws.Range("L6").Value = 2316
$ws.Range("M6").Value = 41.99999800000001
$ws.Range("N6").Value = -2540

# Row 10 (G=1959)
$ws.Range("H10").Value = 0
$ws.Range("I10").Value = 0
$ws.Range("J10").Value = 0
$ws.Range("K10").Value = 0
$ws.Range("L10").Value = 0
$ws.Range("M10").ClearContents()
$ws.Range("N10").ClearContents()

# Row 42 (G=4600)
$ws.Range("H42").Value = 374.16666
$ws.Range("J42").Value = 225.66667
$ws.Range("L42").Value = 677.00001
$ws.Range("N42").Value = -1137.00001

# Row 86 (G=12603)
$ws.Range("H86").Value = 1000
$ws.Range("J86").Value = 1000
$ws.Range("L86").Value = 1000
$ws.Range("N86").Value = -3246

# Row 89 (G=12603)
$ws.Range("H89").Value = 1000
$ws.Range("J89").Value = 1000
$ws.Range("L89").Value = 5000
$ws.Range("N89").Value = -16232

# Row 92 (G=19901)
$ws.Range("H92").Value = 1021.4
$ws.Range("I92").Value = 1021.4
$ws.Range("K92").Value = 1021.4
$ws.Range("M92").Value = 226.6

# Row 137 (G=44013)
$ws.Range("H137").Value = 4381.6665
$ws.Range("I137").Value = 3698.75
$ws.Range("K137").Value = 11096.25
$ws.Range("M137").Value = -8546.25

$ws = $wb.Worksheets.Item("ARM")
# Row 32 (G=44147)
$ws.Range("H32").Value = 8388.817999999999
$ws.Range("I32").Value = 4091.4707
$ws.Range("J32").Value = 22999.8
$ws.Range("K32").Value = 4091.4707
$ws.Range("L32").Value = 22999.8
$ws.Range("M32").Value = -3804.4707
$ws.Range("N32").Value = -23573.8

# Row 45 (G=27714)
$ws.Range("H45").Value = 2000
$ws.Range("I45").Value = 2000
$ws.Range("K45").Value = 2000
$ws.Range("M45").Value = -1623

# Row 74 (G=44000)
$ws.Range("H74").Value = 4083
$ws.Range("I74").Value = 666
$ws.Range("J74").Value = 7500
$ws.Range("K74").Value = 666
$ws.Range("L74").Value = 7500
$ws.Range("M74").Value = 208
$ws.Range("N74").Value = -9248

# Row 77 (G=44000)
$ws.Range("H77").Value = 4083
$ws.Range("I77").Value = 666
$ws.Range("J77").Value = 7500
$ws.Range("K77").Value = 3330
$ws.Range("L77").Value = 37500
$ws.Range("M77").Value = 1038
$ws.Range("N77").Value = -46236

$ws = $wb.Worksheets.Item("BSM")
# Row 25 (G=2370)
$ws.Range("H25").Value = 44599.6
$ws.Range("J25").Value = 44599.6
$ws.Range("L25").Value = 44599.6
$ws.Range("N25").Value = -45069.6

$ws = $wb.Worksheets.Item("CRP")
# Row 11 (G=1821)
$ws.Range("H11").Value = 13332.25
$ws.Range("I11").Value = 20000
$ws.Range("J11").Value = 11109.667
$ws.Range("K11").Value = 20000
$ws.Range("L11").Value = 11109.667
$ws.Range("M11").Value = -19860
$ws.Range("N11").Value = -11389.667

# Row 35 (G=1627)
$ws.Range("H35").Value = 6812.25
$ws.Range("I35").Value = 2282.6667
$ws.Range("J35").Value = 9530
$ws.Range("K35").Value = 2282.6667
$ws.Range("L35").Value = 9530
$ws.Range("M35").Value = -1988.6667
$ws.Range("N35").Value = -10118

# Row 86 (G=12584)
$ws.Range("H86").Value = 1372.5
$ws.Range("I86").Value = 995
$ws.Range("J86").Value = 1750
$ws.Range("K86").Value = 995
$ws.Range("L86").Value = 1750
$ws.Range("M86").Value = 128
$ws.Range("N86").Value = -3996

# Row 89 (G=12584)
$ws.Range("H89").Value = 1372.5
$ws.Range("I89").Value = 995
$ws.Range("J89").Value = 1750
$ws.Range("K89").Value = 4975
$ws.Range("L89").Value = 8750
$ws.Range("M89").Value = 641
$ws.Range("N89").Value = -19982

# Row 132 (G=44019)
$ws.Range("H132").Value = 7381
$ws.Range("I132").Value = 3841.6667
$ws.Range("J132").Value = 17999
$ws.Range("K132").Value = 11525.0001
$ws.Range("L132").Value = 53997
$ws.Range("M132").Value = -8995.000100000001
$ws.Range("N132").Value = -59057

# Row 134 (G=44020)
$ws.Range("H134").Value = 2199.2856
$ws.Range("I134").Value = 2199.2856
$ws.Range("K134").Value = 6597.8568
$ws.Range("M134").Value = -4062.8568

$ws = $wb.Worksheets.Item("CUL")
# Row 33 (G=4867)
$ws.Range("H33").Value = 270
$ws.Range("I33").Value = 390
$ws.Range("J33").Value = 150
$ws.Range("K33").Value = 2340
$ws.Range("L33").Value = 900
$ws.Range("M33").Value = -2057
$ws.Range("N33").Value = -1466

# Row 137 (G=44088)
$ws.Range("H137").Value = 0
$ws.Range("J137").Value = 0
$ws.Range("L137").Value = 0
$ws.Range("N137").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 93 (G=18107)
$ws.Range("H93").Value = 30100.5
$ws.Range("J93").Value = 30100.5
$ws.Range("L93").Value = 30100.5
$ws.Range("N93").Value = -33844.5

# Row 102 (G=36169)
$ws.Range("H102").Value = 416
$ws.Range("I102").Value = 320.125
$ws.Range("K102").Value = 320.125
$ws.Range("M102").Value = 1301.875

$ws = $wb.Worksheets.Item("LTW")
# Row 13 (G=3546)
$ws.Range("H13").Value = 19000
$ws.Range("J13").Value = 19000
$ws.Range("L13").Value = 19000
$ws.Range("N13").Value = -19280

# Row 55 (G=5284)
$ws.Range("H55").Value = 1046.3334
$ws.Range("I55").Value = 961.625
$ws.Range("J55").Value = 1143.1428
$ws.Range("K55").Value = 961.625
$ws.Range("L55").Value = 1143.1428
$ws.Range("M55").Value = -788.625
$ws.Range("N55").Value = -1489.1428

# Row 132 (G=44058)
$ws.Range("H132").Value = 6248
$ws.Range("I132").Value = 5829.3335
$ws.Range("K132").Value = 17488.0005
$ws.Range("M132").Value = -14958.0005

# Row 136 (G=44060)
$ws.Range("H136").Value = 41226
$ws.Range("I136").Value = 41226
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 123678
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -121128
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
# Row 45 (G=21726)
$ws.Range("H45").Value = 28297.75
$ws.Range("J45").Value = 30207
$ws.Range("L45").Value = 30207
$ws.Range("N45").Value = -31189

# Row 62 (G=12589)
$ws.Range("H62").Value = 0
$ws.Range("I62").Value = 0
$ws.Range("K62").Value = 0
$ws.Range("M62").ClearContents()

# Row 65 (G=12589)
$ws.Range("H65").Value = 0
$ws.Range("I65").Value = 0
$ws.Range("K65").Value = 0
$ws.Range("M65").ClearContents()

# Row 132 (G=44029)
$ws.Range("H132").Value = 1241.4615
$ws.Range("I132").Value = 1094.5454
$ws.Range("J132").Value = 2049.5
$ws.Range("K132").Value = 3283.6362
$ws.Range("L132").Value = 6148.5
$ws.Range("M132").Value = -753.6361999999999
$ws.Range("N132").Value = -11208.5
